# Update vm_pu.xlsx results for "case with 380 kV" run.
# Bus voltage magnitudes (p.u.) for rows 2-25 (bus indices 0-23), columns B-F and I-N.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020428732269536
$ws.Range("D2").Value = 1.029579878360528
$ws.Range("E2").Value = 1.030703627554892
$ws.Range("F2").Value = 1.040009003384796
$ws.Range("I2").Value = 1.031838522713
$ws.Range("J2").Value = 1.02562589572549
$ws.Range("K2").Value = 1.032393378082195
$ws.Range("L2").Value = 1.033513870259477
$ws.Range("M2").Value = 1.042792569557691
$ws.Range("N2").Value = 1.012586602164159
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021427042247619
$ws.Range("D3").Value = 1.030009215281493
$ws.Range("E3").Value = 1.031604445166615
$ws.Range("F3").Value = 1.041033913186199
$ws.Range("I3").Value = 1.031917995009093
$ws.Range("J3").Value = 1.026261134385619
$ws.Range("K3").Value = 1.032631910205245
$ws.Range("L3").Value = 1.034222850082924
$ws.Range("M3").Value = 1.043627249909035
$ws.Range("N3").Value = 1.012800030739152
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022073195801218
$ws.Range("D4").Value = 1.030284548244446
$ws.Range("E4").Value = 1.032187420117054
$ws.Range("F4").Value = 1.041696087706062
$ws.Range("I4").Value = 1.031966227086019
$ws.Range("J4").Value = 1.026671796055891
$ws.Range("K4").Value = 1.032783053183499
$ws.Range("L4").Value = 1.034681064828934
$ws.Range("M4").Value = 1.044165732416644
$ws.Range("N4").Value = 1.012937915722453
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02234488156195
$ws.Range("D5").Value = 1.030399702505215
$ws.Range("E5").Value = 1.032432522293346
$ws.Range("F5").Value = 1.041974222520229
$ws.Range("I5").Value = 1.031985737969248
$ws.Range("J5").Value = 1.026844346589642
$ws.Range("K5").Value = 1.032845824479435
$ws.Range("L5").Value = 1.034873567524955
$ws.Range("M5").Value = 1.044391723163206
$ws.Range("N5").Value = 1.012995830222074
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022390501296687
$ws.Range("D6").Value = 1.030419002397035
$ws.Range("E6").Value = 1.032473677141283
$ws.Range("F6").Value = 1.042020908267508
$ws.Range("I6").Value = 1.031988968968795
$ws.Range("J6").Value = 1.026873313192246
$ws.Range("K6").Value = 1.032856318879839
$ws.Range("L6").Value = 1.034905881879459
$ws.Range("M6").Value = 1.044429645248006
$ws.Range("N6").Value = 1.013005551244098
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022076825911914
$ws.Range("D7").Value = 1.030286089285701
$ws.Range("E7").Value = 1.032190695109018
$ws.Range("F7").Value = 1.041699805115163
$ws.Range("I7").Value = 1.031966490803091
$ws.Range("J7").Value = 1.026674102045426
$ws.Range("K7").Value = 1.032783894962176
$ws.Range("L7").Value = 1.03468363757444
$ws.Range("M7").Value = 1.044168753639902
$ws.Range("N7").Value = 1.012938689784859
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020766078276846
$ws.Range("D8").Value = 1.029725485954004
$ws.Range("E8").Value = 1.031008044731515
$ws.Range("F8").Value = 1.040355584882432
$ws.Range("I8").Value = 1.031866040591318
$ws.Range("J8").Value = 1.025840655970959
$ws.Range("K8").Value = 1.032474652992962
$ws.Range("L8").Value = 1.033753585136432
$ws.Range("M8").Value = 1.043074986467198
$ws.Range("N8").Value = 1.012658776246547
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.018457768087167
$ws.Range("D9").Value = 1.028718782224605
$ws.Range("E9").Value = 1.028924763710377
$ws.Range("F9").Value = 1.037979216295448
$ws.Range("I9").Value = 1.031664664384542
$ws.Range("J9").Value = 1.024369131349861
$ws.Range("K9").Value = 1.031905296663936
$ws.Range("L9").Value = 1.032110595887492
$ws.Range("M9").Value = 1.041135351104114
$ws.Range("N9").Value = 1.012163877582756
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.016919847907616
$ws.Range("D10").Value = 1.028035130888171
$ws.Range("E10").Value = 1.027536435949502
$ws.Range("F10").Value = 1.03638988195968
$ws.Range("I10").Value = 1.031514127252027
$ws.Range("J10").Value = 1.023386206876639
$ws.Range("K10").Value = 1.031509430476232
$ws.Range("L10").Value = 1.031012548909293
$ws.Range("M10").Value = 1.039834088109901
$ws.Range("N10").Value = 1.011832846431891
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016254139789963
$ws.Range("D11").Value = 1.027736165877125
$ws.Range("E11").Value = 1.026935410286551
$ws.Range("F11").Value = 1.0357004901356
$ws.Range("I11").Value = 1.031445100098686
$ws.Range("J11").Value = 1.022960143232126
$ws.Range("K11").Value = 1.03133417764659
$ws.Range("L11").Value = 1.030536446309741
$ws.Range("M11").Value = 1.039268704379107
$ws.Range("N11").Value = 1.011689248426187
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016006899484464
$ws.Range("D12").Value = 1.027624677826478
$ws.Range("E12").Value = 1.026712182811233
$ws.Range("F12").Value = 1.035444240108314
$ws.Range("I12").Value = 1.031418884407718
$ws.Range("J12").Value = 1.02280181689338
$ws.Range("K12").Value = 1.031268506039808
$ws.Range("L12").Value = 1.030359504948642
$ws.Range("M12").Value = 1.039058407134772
$ws.Range("N12").Value = 1.011635870963453
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.0160599318602
$ws.Range("D13").Value = 1.027648612224232
$ws.Range("E13").Value = 1.026760064935392
$ws.Range("F13").Value = 1.035499214702692
$ws.Range("I13").Value = 1.031424533800593
$ws.Range("J13").Value = 1.022835781470954
$ws.Range("K13").Value = 1.03128261880584
$ws.Range("L13").Value = 1.030397463789717
$ws.Range("M13").Value = 1.039103529647128
$ws.Range("N13").Value = 1.011647322361776
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016233702148966
$ws.Range("D14").Value = 1.027726959184115
$ws.Range("E14").Value = 1.026916957816993
$ws.Range("F14").Value = 1.035679312078482
$ws.Range("I14").Value = 1.031442944841697
$ws.Range("J14").Value = 1.022947057298798
$ws.Range("K14").Value = 1.031328760927354
$ws.Range("L14").Value = 1.030521822229335
$ws.Range("M14").Value = 1.039251327023179
$ws.Range("N14").Value = 1.011684837018261
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016340772177114
$ws.Range("D15").Value = 1.027775173215705
$ws.Range("E15").Value = 1.027013627401149
$ws.Range("F15").Value = 1.035790252298469
$ws.Range("I15").Value = 1.031454212219877
$ws.Range("J15").Value = 1.023015609098862
$ws.Range("K15").Value = 1.031357114491056
$ws.Range("L15").Value = 1.030598430911882
$ws.Range("M15").Value = 1.039342351662113
$ws.Range("N15").Value = 1.011707945904569
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016964033434091
$ws.Range("D16").Value = 1.028054910493856
$ws.Range("E16").Value = 1.027576326836427
$ws.Range("F16").Value = 1.03643560943833
$ws.Range("I16").Value = 1.031518627507879
$ws.Range("J16").Value = 1.023414473853339
$ws.Range("K16").Value = 1.031520980696898
$ws.Range("L16").Value = 1.031044132821137
$ws.Range("M16").Value = 1.039871570238452
$ws.Range("N16").Value = 1.011842371111255
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017355048091178
$ws.Range("D17").Value = 1.028229597094168
$ws.Range("E17").Value = 1.027929328526753
$ws.Range("F17").Value = 1.036840104490182
$ws.Range("I17").Value = 1.031558005568083
$ws.Range("J17").Value = 1.023664550758411
$ws.Range("K17").Value = 1.031622742976221
$ws.Range("L17").Value = 1.031323538699696
$ws.Range("M17").Value = 1.040203019639469
$ws.Range("N17").Value = 1.011926623198248
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.017583141690467
$ws.Range("D18").Value = 1.028331205053017
$ws.Range("E18").Value = 1.028135240876664
$ws.Range("F18").Value = 1.037075923803407
$ws.Range("I18").Value = 1.031580603174813
$ws.Range("J18").Value = 1.02381037285804
$ws.Range("K18").Value = 1.031681728539209
$ws.Range("L18").Value = 1.031486449530693
$ws.Range("M18").Value = 1.040396162342565
$ws.Range("N18").Value = 1.011975740980657
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017660919335405
$ws.Range("D19").Value = 1.028365802530607
$ws.Range("E19").Value = 1.028205453752356
$ws.Range("F19").Value = 1.037156312420901
$ws.Range("I19").Value = 1.031588245408317
$ws.Range("J19").Value = 1.023860087029973
$ws.Range("K19").Value = 1.031701778159824
$ws.Range("L19").Value = 1.031541987396352
$ws.Range("M19").Value = 1.040461987407649
$ws.Range("N19").Value = 1.011992484622944
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.01731309370851
$ws.Range("D20").Value = 1.028210884209642
$ws.Range("E20").Value = 1.027891453494064
$ws.Range("F20").Value = 1.036796717970084
$ws.Range("I20").Value = 1.03155381903279
$ws.Range("J20").Value = 1.023637724372456
$ws.Range("K20").Value = 1.031611863170149
$ws.Range("L20").Value = 1.031293567490013
$ws.Range("M20").Value = 1.040167477498758
$ws.Range("N20").Value = 1.011917586334673
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016182530209552
$ws.Range("D21").Value = 1.027703900059749
$ws.Range("E21").Value = 1.026870756209622
$ws.Range("F21").Value = 1.035626282829194
$ws.Range("I21").Value = 1.031437539135523
$ws.Range("J21").Value = 1.02291429119643
$ws.Range("K21").Value = 1.031315189076682
$ws.Range("L21").Value = 1.030485204403914
$ws.Range("M21").Value = 1.039207812336121
$ws.Range("N21").Value = 1.011673790955667
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.0154718929247
$ws.Range("D22").Value = 1.027382599495569
$ws.Range("E22").Value = 1.026229120651583
$ws.Range("F22").Value = 1.034889347431826
$ws.Range("I22").Value = 1.031361097902842
$ws.Range("J22").Value = 1.022459050683295
$ws.Range("K22").Value = 1.031125333056326
$ws.Range("L22").Value = 1.029976401312478
$ws.Range("M22").Value = 1.038602763822839
$ws.Range("N22").Value = 1.011520282907029
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.015848596897861
$ws.Range("D23").Value = 1.027553166905001
$ws.Range("E23").Value = 1.026569252425848
$ws.Range("F23").Value = 1.035280108736269
$ws.Range("I23").Value = 1.031401936131991
$ws.Range("J23").Value = 1.02270041906748
$ws.Range("K23").Value = 1.031226293854316
$ws.Range("L23").Value = 1.030246179677592
$ws.Range("M23").Value = 1.038923669309951
$ws.Range("N23").Value = 1.011601681618338
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017332051021947
$ws.Range("D24").Value = 1.028219340633712
$ws.Range("E24").Value = 1.027908567552532
$ws.Range("F24").Value = 1.036816322828647
$ws.Range("I24").Value = 1.031555711894149
$ws.Range("J24").Value = 1.023649846195956
$ws.Range("K24").Value = 1.031616780431601
$ws.Range("L24").Value = 1.031307110380728
$ws.Range("M24").Value = 1.040183538037239
$ws.Range("N24").Value = 1.011921669781672
$ws.Range("B25").Value = 1.019999999999999
$ws.Range("C25").Value = 1.019054354681104
$ws.Range("D25").Value = 1.028981256504333
$ws.Range("E25").Value = 1.029463253028167
$ws.Range("F25").Value = 1.038594466824935
$ws.Range("I25").Value = 1.031719602601554
$ws.Range("J25").Value = 1.024749894367964
$ws.Range("K25").Value = 1.032055370493441
$ws.Range("L25").Value = 1.032535830985353
$ws.Range("M25").Value = 1.04163823971371
$ws.Range("N25").Value = 1.012292015425343

$wb.Save()
